# Reorganize the "Alchemist" / "Mage" sections of the sound list:
#  - Insert 4 new "Potion_Choice_0N" sound rows right after the Alchemist
#    header row, pushing the existing Alchemist potion rows (Healing,
#    Poison, Mana, Antidote) down by 4 rows.
#  - The existing Mage "Orb_Attack_*"/"Fireball"/... rows likewise get
#    pushed down so they sit right before the trailing blank Mage rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the old Alchemist potion rows (43-46) before we overwrite them ---
$oldPotionB = @()
$oldPotionC = @()
for ($r = 43; $r -le 46; $r++) {
    $oldPotionB += , $ws.Cells.Item($r, 2).Value()
    $oldPotionC += , $ws.Cells.Item($r, 3).Value()
}

# --- capture the old Mage attack rows (47-55) before we overwrite them ---
$oldMageB = @()
$oldMageC = @()
for ($r = 47; $r -le 55; $r++) {
    $oldMageB += , $ws.Cells.Item($r, 2).Value()
    $oldMageC += , $ws.Cells.Item($r, 3).Value()
}

# --- rows 43-46: brand-new "Potion_Choice" entries ---
$newNames = @("Potion_Choice_01", "Potion_Choice_02", "Potion_Choice_03", "Potion_Choice_04")
$newTimes = @(0.4, 0.37, 0.33, 0.39)
for ($i = 0; $i -lt 4; $i++) {
    $r = 43 + $i
    $ws.Cells.Item($r, 1).Value = "Alchemist"
    $ws.Cells.Item($r, 2).Value = $newNames[$i]
    $ws.Cells.Item($r, 3).Value = $newTimes[$i]
    $ws.Cells.Item($r, 4).Value = 20
}

# --- rows 47-50: the old Alchemist potion rows, shifted down by 4 ---
# (also carries the Alchemist row shading with it, matching row 43-46's look)
for ($i = 0; $i -lt 4; $i++) {
    $r = 47 + $i
    $ws.Cells.Item($r, 1).Value = "Alchemist"
    $ws.Cells.Item($r, 2).Value = $oldPotionB[$i]
    $ws.Cells.Item($r, 3).Value = $oldPotionC[$i]
    $ws.Cells.Item($r, 4).Value = 20

    $ws.Range("A43:D43").Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- rows 51-59: now-blank Mage filler rows ---
for ($r = 51; $r -le 59; $r++) {
    $ws.Cells.Item($r, 1).Value = "Mage"
    $ws.Cells.Item($r, 2).Value = ""
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
}

# --- rows 60-68: the old Mage attack rows, shifted down by 13 ---
for ($i = 0; $i -lt 9; $i++) {
    $r = 60 + $i
    $ws.Cells.Item($r, 1).Value = "Mage"
    $ws.Cells.Item($r, 2).Value = $oldMageB[$i]
    $ws.Cells.Item($r, 3).Value = $oldMageC[$i]
    $ws.Cells.Item($r, 4).Value = 20
}

# --- rows 69-71 stay blank Mage filler rows ---
for ($r = 69; $r -le 71; $r++) {
    $ws.Cells.Item($r, 1).Value = "Mage"
    $ws.Cells.Item($r, 2).Value = ""
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
}

$ws.Range("C47").Select()
